$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as text strings in the
# source data (e.g. "307.12", "2.44%"), not numbers. Force text format
# before assignment so Excel does not auto-convert them to numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.44%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.098"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.38%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08154"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.96%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.935"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.41%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.195"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.24%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.776"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.56%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9294"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.84%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1379"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "25.41%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1927"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.34%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09183"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.31%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03606"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.14%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09846"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.27%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001414"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.98%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005886"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.31%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.594"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.99%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.982"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.88%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3439"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1322"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.10%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.23%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2409"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.42%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04516"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.41%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.19%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004880"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.42%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001242"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.60%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02014"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.75%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04940"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.58%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "16.61%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007640"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.46%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1381"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.51%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.74%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01069"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.49%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006447"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.24%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.07%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001191"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-8.68%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.07%"
